$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.636.92"
$ws.Range("E2").Value = "  +0.97%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.868.04"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.42"
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4712"
$ws.Range("E7").Value = "  -1.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2777"
$ws.Range("E8").Value = "  +0.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06387"
$ws.Range("E9").Value = "  -0.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.06"
$ws.Range("E10").Value = "  +11.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.853.94"
$ws.Range("E11").Value = "  -0.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07463"
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.989"
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("E14").Value = "  -0.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6377"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.598.20"
$ws.Range("E16").Value = "  +1.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "241.94"
$ws.Range("E17").Value = "  +4.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.87"
$ws.Range("E19").Value = "  +0.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007387"
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9998"
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.011"
$ws.Range("E22").Value = "  -1.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.030"
$ws.Range("E23").Value = "  +0.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.398"
$ws.Range("E24").Value = "  +1.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "166.19"
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.28"
$ws.Range("E26").Value = "  +2.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.895"
$ws.Range("E27").Value = "  +1.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1028"
$ws.Range("E28").Value = "  +2.64%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.110"
$ws.Range("E30").Value = "  -2.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.878"
$ws.Range("E31").Value = "  -1.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04940"
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.159"
$ws.Range("E33").Value = "  +1.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7111"
$ws.Range("E34").Value = "  -1.64%  "
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.01909"
$ws.Range("E36").Value = "  -0.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.697"
$ws.Range("E37").Value = "  +2.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.8836"
$ws.Range("E38").Value = "  -2.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.002"
$ws.Range("E39").Value = "  +0.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "106.00"
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4125"
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.555"
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.509"
$ws.Range("E44").Value = "  +6.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.46"
$ws.Range("E45").Value = "  +1.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1230"
$ws.Range("E46").Value = "  +1.89%  "
$ws.Range("B47").Value = "Elrond"
$ws.Range("C47").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "33.74"
$ws.Range("E47").Value = "  +2.01%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.624"
$ws.Range("E48").Value = "  -1.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.390"
$ws.Range("E49").Value = "  -0.94%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05580"
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3728"
$ws.Range("E51").Value = "  +0.63%  "
